# Add an "Example Grid" example/reference sheet next to the existing
# start-cell calculation sheet, and give that first sheet a clearer name.

$wb = $excel.ActiveWorkbook

# --- Rename the original sheet ---------------------------------------
$startCells = $wb.Worksheets.Item(1)
$startCells.Name = "Start Cells"

# --- Add the new "Example Grid" sheet right after it ------------------
$grid = $wb.Worksheets.Add($null, $startCells)
$grid.Name = "Example Grid"

# 5x5 grid of sequential numbers (0-24), row by row, starting at A1
$gridValues = @(
    @(0, 1, 2, 3, 4),
    @(5, 6, 7, 8, 9),
    @(10, 11, 12, 13, 14),
    @(15, 16, 17, 18, 19),
    @(20, 21, 22, 23, 24)
)

for ($r = 0; $r -lt 5; $r++) {
    for ($c = 0; $c -lt 5; $c++) {
        $grid.Cells.Item($r + 1, $c + 1).Value = $gridValues[$r][$c]
    }
}

# Labelled helper value off to the side of the grid
$grid.Range("L1").Value = "grid_size"
$grid.Range("L1").Font.Bold = $true
$grid.Range("L2").Value = 5

# Widen the label column so "grid_size" fits comfortably
$grid.Columns.Item(12).ColumnWidth = 13.35

$grid.PageSetup.Orientation = 1
